$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6815946698188782
$ws.Range("B1").Value = 1.525940299034119
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.020338296890259
$ws.Range("E1").Value = 1.268260478973389
